$d = $word.ActiveDocument

# 1 & 2: Mark the two inline pictures' ranges as "no proofing" so Word
# writes <w:noProof/> into the rPr of the runs that host the <w:drawing>
# elements (matches author re-saving the doc after the pictures were
# (re)inserted/viewed).
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}

# 3: Add two new paragraphs right after the paragraph ending in
# "...found in the PGN." -- first an empty paragraph, then one with the
# footnote-style clarification text.
$find = $d.Content
$found = $find.Find.Execute("30-character string into sets of 3, and encode it using number found in the PGN.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertAt = $find.End
    $target = $d.Range($insertAt, $insertAt)
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>*Sets of 3 because chess notation is in sets of 3</w:t></w:r></w:p>'
    $target.InsertXML($xml)
}
